$d = $word.ActiveDocument

# --- Locate the block of paragraphs to replace -----------------------------
# The block starts at the paragraph "Fourth Video: Menu Overlay & Responsiveness"
# and ends at the empty paragraph right after "Media Query Mixins in
# _config.scss, created _mobile.scss" (the paragraph right before the
# "EXAMPLE STYLE 2" page-break heading).

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -clike "Fourth Video: Menu Overlay*") {
        $startIdx = $i
    }
    if ($t -clike "*Media Query Mixins*config.scss*") {
        # the block ends with the (empty) paragraph right after this one
        $endIdx = $i + 1
    }
}

$startPara = $d.Paragraphs.Item($startIdx)
$endPara = $d.Paragraphs.Item($endIdx)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newXml = @'
<w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I started watching the fourth video about the Menu Overlay &amp; Responsiveness. I’ve learned how to use the &amp;-character</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to add classes to the class I am already in in the Sass-file. I also learned how to position specific items vertically and horizontally</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and how to style them</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>. To avoid scroll bars, overflow always needs to be hidden. To darken or lighten colors a little bit, there are useful functions that need the color and the parameter of how much one wants to darken/lighten that specific color.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I already knew about list-style that hides the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bulletpoints</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> if it is set to none.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>An important function that was shown in the video was translate3</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>d(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>). It helps for example the menu to slide into the viewport smoothly, depending on how the parameters are set (if from the top, left, right, or bottom).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I also learned how to implement a delay for each individual item</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>nth-child</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, which can also be useful in the future.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>27.12.2022</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">I continued working on the fourth video and learned a bit about responsiveness. From what I’ve learned in the past I already knew about the four different screen sizes that are also being used in the video. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Still,</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> I had no clue about how to implement them as Sass functions </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Mixins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">as I </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>did</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> not work with Sass before.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I think the responsive part of the example project would be best to use for my own project as well, because I cannot think of any reason why I would need to change that, if it works so well.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Also, this might be my favorite part of the course until now, because responsiveness is such a useful topic to go through and it is easy to understand why the person in the video is doing what he is doing and how one could change that according to their own wishes.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>28.12.2022</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I worked through the fifth video and therefore</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> created a color function which I will probably be using in my own project as well. It makes sure that there won’t be any unreadable texts in terms of how background and text color work together.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Textkrper"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The “About Me”-html page part was not very new to me, because I have worked with html a little before, but what I learned was another way of creating the CSS grid. </w:t></w:r></w:p>
'@

$r.InsertXML($newXml)
